$d = $word.ActiveDocument

# ------------------------------------------------------------------
# The "Models" table (Table 1) lists each model column together with
# a "filter/search" flag. Two edits land in that table:
#   1. The row for "secondLanguage" keeps its "filter" text, but now
#      carries the (previously end-of-document) "_GoBack" bookmark
#      right after that run.
#   2. The row for "city" changes its flag from "filter" to "search"
#      and is now bold (cities are searchable).
# ------------------------------------------------------------------
$tbl = $d.Tables.Item(1)

# Locate rows by their first-column label so this keeps working even
# if row ordering ever shifts.
$secondLanguageRow = 0
$cityRow = 0
for ($i = 1; $i -le $tbl.Rows.Count; $i++) {
    $label = $tbl.Cell($i, 1).Range.Text
    if ($label -like "secondLanguage*") { $secondLanguageRow = $i }
    if ($label -like "city*") { $cityRow = $i }
}

# --- Edit 1: change "city" row's filter/search cell to bold "search" ---
$cityCell = $tbl.Cell($cityRow, 3)
$cr = $cityCell.Range
$cityWordRange = $d.Range($cr.Start, $cr.Start + 6)
$cityWordRange.Text = "search"
$cityWordRange2 = $d.Range($cr.Start, $cr.Start + 6)
$cityWordRange2.Bold = 1
# Also bold the paragraph mark itself, matching a whole-paragraph bold.
$cityPara = $cityCell.Range.Paragraphs.Item(1)
$cityPara.Range.Bold = 1

# --- Edit 2: move the "_GoBack" bookmark onto the "secondLanguage" row's
#     filter run (collapsed, right after the word). ---
$slCell = $tbl.Cell($secondLanguageRow, 3)
$slRange = $slCell.Range
$slWordRange = $d.Range($slRange.Start, $slRange.Start + 6)
$slEnd = $slWordRange.End

# A zero-length Range positioned exactly on the trailing paragraph/cell
# mark isn't accepted directly, so nudge it into existence by inserting
# a throwaway character, anchoring the bookmark there, then removing
# the throwaway character again (the bookmark stays put, collapsed).
$anchor = $d.Range($slEnd, $slEnd)
$anchor.InsertAfter("X")
$bmRange = $d.Range($slEnd, $slEnd)
$d.Bookmarks.Add("_GoBack", $bmRange)
$throwaway = $d.Range($slEnd, $slEnd + 1)
$throwaway.Delete()
